$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 341; this shifts the existing rows 341-360
# down to 342-361 (and carries the D-column date style down with them).
$ws.Rows.Item(341).EntireRow.Insert()

# Populate the newly inserted row 341 with the new weekly record.
$ws.Range("A341").Value = 3
$ws.Range("B341").Value = "Femacal de La Calera"
$ws.Range("C341").Value = "Coquimbo"
$ws.Range("D341").Value = 44706
$ws.Range("E341").Value = 5
$ws.Range("F341").Value = 100112040
$ws.Range("G341").Value = "Cilantro"
$ws.Range("H341").Value = "Sin especificar"
$ws.Range("I341").Value = "Primera"
$ws.Range("J341").Value = 170
$ws.Range("K341").Value = 3300
$ws.Range("L341").Value = 3500
$ws.Range("M341").Value = 3406
$ws.Range("N341").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O341").Value = "Provincia de Quillota"
$ws.Range("P341").Value = 1135
$ws.Range("Q341").Value = 3
$ws.Range("R341").Value = "Hortaliza"
